$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Minimum_daily_mileage"
$ws.Range("A8").Value = "Battery_capacity"
$ws.Range("A9").Value = "Electric_consumption_NEFZ"
$ws.Range("A10").Value = "Fuel_consumption_NEFZ"
$ws.Range("A11").Value = "Electric_consumption_Artemis"
$ws.Range("A12").Value = "Fuel_consumption_Artemis"
$ws.Range("A13").Value = "Maximum_SOC"
$ws.Range("A14").Value = "Minimum_SOC"
$ws.Range("A15").Value = "Rated_power_of_charging_column"
$ws.Range("A16").Value = "Is_BEV?"

$ws.Range("C21").Select()
